$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cell, $value) {
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.Style = "Normal"
}

Set-TextValue $ws.Range('D2') '91.829.77'
Set-TextValue $ws.Range('E2') '  +1.33%  '
Set-TextValue $ws.Range('D3') '3.129.00'
Set-TextValue $ws.Range('E3') '  -0.45%  '
Set-TextValue $ws.Range('D4') '1.00'
Set-TextValue $ws.Range('E4') '  +0.02%  '
Set-TextValue $ws.Range('D5') '241.87'
Set-TextValue $ws.Range('E5') '  +1.62%  '
Set-TextValue $ws.Range('D6') '626.68'
Set-TextValue $ws.Range('E6') '  -2.33%  '
Set-TextValue $ws.Range('D7') '1.17'
Set-TextValue $ws.Range('E7') '  +7.77%  '
Set-TextValue $ws.Range('D8') '0.375'
Set-TextValue $ws.Range('E8') '  +4.38%  '
Set-TextValue $ws.Range('D9') '1.00'
Set-TextValue $ws.Range('E9') '  -0.01%  '
Set-TextValue $ws.Range('D10') '3.128.09'
Set-TextValue $ws.Range('E10') '  -0.41%  '
Set-TextValue $ws.Range('D11') '0.764'
Set-TextValue $ws.Range('E11') '  +5.56%  '
Set-TextValue $ws.Range('D12') '0.205'
Set-TextValue $ws.Range('E12') '  +3.72%  '
Set-TextValue $ws.Range('D13') '0.0000254'
Set-TextValue $ws.Range('E13') '  +3.44%  '
Set-TextValue $ws.Range('D14') '35.86'
Set-TextValue $ws.Range('E14') '  -2.47%  '
Set-TextValue $ws.Range('E15') '  -2.10%  '
Set-TextValue $ws.Range('D16') '91.325.96'
Set-TextValue $ws.Range('E16') '  +1.12%  '
Set-TextValue $ws.Range('E17') '  -0.24%  '
Set-TextValue $ws.Range('D18') '3.143.85'
Set-TextValue $ws.Range('E18') '  +2.23%  '
Set-TextValue $ws.Range('D19') '3.80'
Set-TextValue $ws.Range('E19') '  +2.17%  '
Set-TextValue $ws.Range('D20') '14.78'
Set-TextValue $ws.Range('E20') '  +1.79%  '
Set-TextValue $ws.Range('E21') '  -1.07%  '
Set-TextValue $ws.Range('D22') '5.88'
Set-TextValue $ws.Range('E22') '  +3.12%  '
Set-TextValue $ws.Range('D23') '452.19'
Set-TextValue $ws.Range('E23') '  +0.14%  '
Set-TextValue $ws.Range('D24') '9.22'
Set-TextValue $ws.Range('E24') '  +1.62%  '
Set-TextValue $ws.Range('E25') '  -1.50%  '
Set-TextValue $ws.Range('D26') '93.51'
Set-TextValue $ws.Range('E26') '  +2.33%  '
Set-TextValue $ws.Range('D27') '12.12'
Set-TextValue $ws.Range('E27') '  -3.25%  '
Set-TextValue $ws.Range('D28') '3.297.27'
Set-TextValue $ws.Range('E28') '  -0.08%  '
Set-TextValue $ws.Range('E29') '  +0.12%  '
Set-TextValue $ws.Range('D30') '0.181'
Set-TextValue $ws.Range('E30') '  +12.97%  '
Set-TextValue $ws.Range('D31') '0.236'
Set-TextValue $ws.Range('E31') '  +16.92%  '
Set-TextValue $ws.Range('D32') '0.118'
Set-TextValue $ws.Range('E32') '  +36.41%  '
Set-TextValue $ws.Range('D33') '9.28'
Set-TextValue $ws.Range('E33') '  -6.85%  '
Set-TextValue $ws.Range('E34') '  +36.33%  '
Set-TextValue $ws.Range('E35') '  +9.30%  '
Set-TextValue $ws.Range('E36') '  -1.73%  '
Set-TextValue $ws.Range('D37') '7.59'
Set-TextValue $ws.Range('E37') '  +6.21%  '
Set-TextValue $ws.Range('D38') '4.16'
Set-TextValue $ws.Range('E38') '  +24.83%  '
Set-TextValue $ws.Range('D39') '500.53'
Set-TextValue $ws.Range('E39') '  -3.62%  '
Set-TextValue $ws.Range('E40') '  -0.84%  '
Set-TextValue $ws.Range('D41') '3.65'
Set-TextValue $ws.Range('E41') '  -6.64%  '
Set-TextValue $ws.Range('E42') '  -0.66%  '
Set-TextValue $ws.Range('E43') '  -0.15%  '
Set-TextValue $ws.Range('D44') '22.18'
Set-TextValue $ws.Range('E44') '  -0.15%  '
Set-TextValue $ws.Range('E46') '  -0.18%  '
Set-TextValue $ws.Range('D47') '157.50'
Set-TextValue $ws.Range('E47') '  +4.96%  '
Set-TextValue $ws.Range('E48') '  -0.67%  '
Set-TextValue $ws.Range('E49') '  -0.39%  '
Set-TextValue $ws.Range('E50') '  -0.23%  '
Set-TextValue $ws.Range('D51') '44.93'
Set-TextValue $ws.Range('E51') '  -1.65%  '
